$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2, 2, 2, 1, 1, 1, 2, 1, 2, 2, 2, 1, 2, 2)
$col = 2
foreach ($v in $values) {
    $ws.Cells.Item(2, $col).Value = $v
    $col = $col + 1
}
